$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the existing "notes" column (H), shifting
# notes (and the styled blank marker cell) one column to the right.
$ws.Columns("H:H").Insert()

# Header for the new column.
$ws.Range("H1").Value = "unique speakers"

# Per-dataset unique-speaker counts (row 6 keeps a text suffix instead of a count).
$ws.Range("H2").Value = 6
$ws.Range("H4").Value = 30
$ws.Range("H5").Value = 147
$ws.Range("H6").Value = " [tr+en]"
$ws.Range("H7").Value = 12
$ws.Range("H8").Value = 91
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 10
$ws.Range("H12").Value = 63
$ws.Range("H13").Value = 3
$ws.Range("H14").Value = 4
$ws.Range("H15").Value = 44
$ws.Range("H16").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("H18").Value = 6
$ws.Range("H19").Value = 4
$ws.Range("H20").Value = 134

# Resize the new column to fit its header text.
$ws.Columns("H:H").AutoFit()

# Update the LEGOv2 notes to mention the unknown gender labels.
$ws.Range("I20").Value = "from an automated bus info service; some unknown gender labels"

# New dataset row: MELD.
$ws.Range("A21").Value = "MELD"

# Move the active selection like the author left it.
$ws.Range("B21").Select()
